$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their literal text representation
# (many values look numeric, e.g. "1.00", "0.621", and would otherwise be
# auto-converted to numbers by Excel, losing formatting / trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "43.627.97"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.284.67"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "112.00"
$ws.Range("E5").Value = "  +17.11%  "
$ws.Range("D6").Value = "267.07"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").Value = "0.612"
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("D10").Value = "47.17"
$ws.Range("E10").Value = "  +5.05%  "
$ws.Range("D11").Value = "0.0933"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "8.49"
$ws.Range("E12").Value = "  +8.69%  "
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").Value = "15.48"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "2.626.32"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "0.846"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "2.285.72"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "43.461.72"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "6.51"
$ws.Range("E20").Value = "  +5.46%  "
$ws.Range("D21").Value = "72.04"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "2.51"
$ws.Range("E22").Value = "  +3.20%  "
$ws.Range("D23").Value = "232.23"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").Value = "9.44"
$ws.Range("E24").Value = "  +4.08%  "
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +12.60%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "11.34"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").Value = "42.90"
$ws.Range("E28").Value = "  +5.99%  "
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "175.72"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").Value = "21.59"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").Value = "0.0922"
$ws.Range("E33").Value = "  +4.49%  "
$ws.Range("D34").Value = "5.46"
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("E36").Value = "  +6.95%  "
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0351"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "3.84"
$ws.Range("E39").Value = "  +15.33%  "
$ws.Range("D40").Value = "74.00"
$ws.Range("E40").Value = "  +16.86%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "2.41"
$ws.Range("E41").Value = "  +3.51%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.241"
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("D43").Value = "13.21"
$ws.Range("E43").Value = "  +8.67%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "1.41"
$ws.Range("E45").Value = "  +5.18%  "
$ws.Range("D46").Value = "5.93"
$ws.Range("E46").Value = "  +13.39%  "
$ws.Range("D47").Value = "8.74"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").Value = "0.0998"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").Value = "101.01"
$ws.Range("E49").Value = "  +2.98%  "
$ws.Range("D50").Value = "1.22"
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("D51").Value = "0.450"
$ws.Range("E51").Value = "  +5.73%  "
